$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two shared-string text values (Cigarrette Other -> Tobacco Other, Empty -> General Empty)
$ws.Range("A22").Value = "Tobacco Other"
$ws.Range("A23").Value = "General Empty"

# Move the sheet's active selection from C1 to A22
$ws.Range("A22").Select()

# Widen columns A and B (and the rest of B's style range) slightly
$ws.Columns.Item(1).ColumnWidth = 30.6703703703704
$ws.Columns.Item(2).ColumnWidth = 12.4444444444444

# Nudge the workbook's tab ratio (horizontal scroll/tab split) by one unit
$wb.Windows.Item(1).TabRatio = 993
